$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 232, shifting existing rows
# 232:261 down to 236:265 (dimension grows from A1:T261 to A1:T265).
$ws.Rows("232:235").Insert()

# --- New row 232: Naranja / Lane Late / Primera ---
$ws.Cells.Item(232, 1).Value = 4
$ws.Cells.Item(232, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(232, 3).Value = "Los Lagos"
$ws.Cells.Item(232, 4).Value = 44491
$ws.Cells.Item(232, 5).Value = 10
$ws.Cells.Item(232, 6).Value = "Fruta"
$ws.Cells.Item(232, 7).Value = 100102
$ws.Cells.Item(232, 8).Value = "Cítricos"
$ws.Cells.Item(232, 9).Value = 100102005
$ws.Cells.Item(232, 10).Value = "Naranja"
$ws.Cells.Item(232, 11).Value = "Lane Late"
$ws.Cells.Item(232, 12).Value = "Primera"
$ws.Cells.Item(232, 13).Value = 300
$ws.Cells.Item(232, 14).Value = 13000
$ws.Cells.Item(232, 15).Value = 13000
$ws.Cells.Item(232, 16).Value = 13000
$ws.Cells.Item(232, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(232, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(232, 19).Value = 867
$ws.Cells.Item(232, 20).Value = 15

# --- New row 233: Naranja / Lane Late / Segunda ---
$ws.Cells.Item(233, 1).Value = 4
$ws.Cells.Item(233, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(233, 3).Value = "Los Lagos"
$ws.Cells.Item(233, 4).Value = 44491
$ws.Cells.Item(233, 5).Value = 10
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100102
$ws.Cells.Item(233, 8).Value = "Cítricos"
$ws.Cells.Item(233, 9).Value = 100102005
$ws.Cells.Item(233, 10).Value = "Naranja"
$ws.Cells.Item(233, 11).Value = "Lane Late"
$ws.Cells.Item(233, 12).Value = "Segunda"
$ws.Cells.Item(233, 13).Value = 300
$ws.Cells.Item(233, 14).Value = 11000
$ws.Cells.Item(233, 15).Value = 11000
$ws.Cells.Item(233, 16).Value = 11000
$ws.Cells.Item(233, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(233, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(233, 19).Value = 733
$ws.Cells.Item(233, 20).Value = 15

# --- New row 234: Naranja / Navel Late / Primera ---
$ws.Cells.Item(234, 1).Value = 4
$ws.Cells.Item(234, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(234, 3).Value = "Los Lagos"
$ws.Cells.Item(234, 4).Value = 44491
$ws.Cells.Item(234, 5).Value = 10
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100102
$ws.Cells.Item(234, 8).Value = "Cítricos"
$ws.Cells.Item(234, 9).Value = 100102005
$ws.Cells.Item(234, 10).Value = "Naranja"
$ws.Cells.Item(234, 11).Value = "Navel Late"
$ws.Cells.Item(234, 12).Value = "Primera"
$ws.Cells.Item(234, 13).Value = 300
$ws.Cells.Item(234, 14).Value = 13000
$ws.Cells.Item(234, 15).Value = 13000
$ws.Cells.Item(234, 16).Value = 13000
$ws.Cells.Item(234, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(234, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(234, 19).Value = 867
$ws.Cells.Item(234, 20).Value = 15

# --- New row 235: Naranja / Navel Late / Segunda ---
$ws.Cells.Item(235, 1).Value = 4
$ws.Cells.Item(235, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(235, 3).Value = "Los Lagos"
$ws.Cells.Item(235, 4).Value = 44491
$ws.Cells.Item(235, 5).Value = 10
$ws.Cells.Item(235, 6).Value = "Fruta"
$ws.Cells.Item(235, 7).Value = 100102
$ws.Cells.Item(235, 8).Value = "Cítricos"
$ws.Cells.Item(235, 9).Value = 100102005
$ws.Cells.Item(235, 10).Value = "Naranja"
$ws.Cells.Item(235, 11).Value = "Navel Late"
$ws.Cells.Item(235, 12).Value = "Segunda"
$ws.Cells.Item(235, 13).Value = 300
$ws.Cells.Item(235, 14).Value = 11000
$ws.Cells.Item(235, 15).Value = 11000
$ws.Cells.Item(235, 16).Value = 11000
$ws.Cells.Item(235, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(235, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(235, 19).Value = 733
$ws.Cells.Item(235, 20).Value = 15
